$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.086.38"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.666.85"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'216.64"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "'0.5112"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.2630"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").Value = "'0.06416"
$ws.Range("E9").Value = "  +4.11%  "
$ws.Range("D10").Value = "'21.71"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "'0.07426"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "1.672.74"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "'4.509"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").Value = "'0.5815"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "'0.000008564"
$ws.Range("E15").Value = "  +5.09%  "
$ws.Range("D16").Value = "'64.37"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").Value = "26.155.54"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "'4.922"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "'10.78"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "'189.01"
$ws.Range("E21").Value = "  +3.19%  "
$ws.Range("D22").Value = "'6.205"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").Value = "'1.006"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'145.86"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Value = "'7.630"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("D26").Value = "'0.1191"
$ws.Range("E26").Value = "  +5.86%  "
$ws.Range("D27").Value = "'15.62"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("D28").Value = "'0.06412"
$ws.Range("E28").Value = "  +13.17%  "
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").Value = "'1.320"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'3.524"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").Value = "'3.512"
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("D33").Value = "'1.638"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").Value = "'0.6073"
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").Value = "'2.367"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "'2.687"
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").Value = "'6.201"
$ws.Range("E38").Value = "  +5.85%  "
$ws.Range("D39").Value = "'0.01613"
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("D40").Value = "1.076.05"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "'0.8616"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").Value = "1.815.51"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("E45").Value = "  +9.32%  "
$ws.Range("D46").Value = "'56.18"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'1.005"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "'8.065"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "'0.05206"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").Value = "'0.4292"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "'5.951"
$ws.Range("E51").Value = "  +6.35%  "
